# Applies the Rafflesia Profits market-data refresh produced by the
# scheduled runner: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) for the affected Leve rows on each job worksheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 22.6
$ws.Range("I8").Value = 22.6
$ws.Range("K8").Value = 67.80000000000001
$ws.Range("M8").Value = 71.19999999999999
$ws.Range("H33").Value = 340.41666
$ws.Range("I33").Value = 287.1111
$ws.Range("K33").Value = 287.1111
$ws.Range("M33").Value = -58.11110000000002
$ws.Range("H80").Value = 1161.7693
$ws.Range("I80").Value = 1204.1428
$ws.Range("J80").Value = 1112.3334
$ws.Range("K80").Value = 3612.4284
$ws.Range("L80").Value = 3337.0002
$ws.Range("M80").Value = -2614.4284
$ws.Range("N80").Value = -5333.0002
$ws.Range("H83").Value = 1161.7693
$ws.Range("I83").Value = 1204.1428
$ws.Range("J83").Value = 1112.3334
$ws.Range("K83").Value = 10837.2852
$ws.Range("L83").Value = 10011.0006
$ws.Range("M83").Value = -5845.2852
$ws.Range("N83").Value = -19995.0006
$ws.Range("H116").Value = 2500
$ws.Range("I116").Value = 2500
$ws.Range("K116").Value = 2500
$ws.Range("M116").Value = 942
$ws.Range("H132").Value = 6198.25
$ws.Range("I132").Value = 6037.9
$ws.Range("K132").Value = 18113.7
$ws.Range("M132").Value = -15583.7
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 425
$ws.Range("I2").Value = 425
$ws.Range("K2").Value = 425
$ws.Range("M2").Value = -312
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = $null
$ws.Range("H12").Value = 226.5
$ws.Range("I12").Value = 226.5
$ws.Range("K12").Value = 226.5
$ws.Range("M12").Value = -53.5
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = $null
$ws.Range("H19").Value = 254
$ws.Range("I19").Value = 254
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 254
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -25
$ws.Range("N19").Value = $null
$ws.Range("H116").Value = 425
$ws.Range("I116").Value = 425
$ws.Range("K116").Value = 425
$ws.Range("M116").Value = 1869

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 425
$ws.Range("I3").Value = 425
$ws.Range("K3").Value = 425
$ws.Range("M3").Value = -311
$ws.Range("H11").Value = 3027.5
$ws.Range("I11").Value = 500
$ws.Range("J11").Value = 5555
$ws.Range("K11").Value = 500
$ws.Range("L11").Value = 5555
$ws.Range("M11").Value = -360
$ws.Range("N11").Value = -5835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = $null
$ws.Range("H4").Value = 77.5
$ws.Range("J4").Value = 77.5
$ws.Range("L4").Value = 77.5
$ws.Range("N4").Value = -301.5
$ws.Range("H19").Value = 342.4
$ws.Range("I19").Value = 199.83333
$ws.Range("J19").Value = 437.44446
$ws.Range("K19").Value = 199.83333
$ws.Range("L19").Value = 437.44446
$ws.Range("M19").Value = -29.83332999999999
$ws.Range("N19").Value = -777.4444599999999
$ws.Range("H24").Value = 342.4
$ws.Range("I24").Value = 199.83333
$ws.Range("J24").Value = 437.44446
$ws.Range("K24").Value = 199.83333
$ws.Range("L24").Value = 437.44446
$ws.Range("M24").Value = -29.83332999999999
$ws.Range("N24").Value = -777.4444599999999
$ws.Range("H88").Value = 29990
$ws.Range("J88").Value = 29990
$ws.Range("L88").Value = 29990
$ws.Range("N88").Value = -30802
$ws.Range("H91").Value = 29990
$ws.Range("J91").Value = 29990
$ws.Range("L91").Value = 29990
$ws.Range("N91").Value = -32798

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 96.888885
$ws.Range("I17").Value = 105.25
$ws.Range("K17").Value = 315.75
$ws.Range("M17").Value = -146.75
$ws.Range("H68").Value = 802
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = $null
$ws.Range("H71").Value = 802
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 636.6667
$ws.Range("I5").Value = 636.6667
$ws.Range("K5").Value = 636.6667
$ws.Range("M5").Value = -524.6667
$ws.Range("H62").Value = 30000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 30000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 30000
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -31372
$ws.Range("H65").Value = 30000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 30000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 90000
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -96864

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = $null
$ws.Range("H9").Value = 1747.6666
$ws.Range("I9").Value = 1747.6666
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1747.6666
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -1523.6666
$ws.Range("N9").Value = $null
$ws.Range("H10").Value = 500
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 500
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 500
$ws.Range("M10").Value = $null
$ws.Range("N10").Value = -780
$ws.Range("H12").Value = 4000
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 4000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 4000
$ws.Range("M12").Value = $null
$ws.Range("N12").Value = -4340
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = $null
$ws.Range("H17").Value = 2999
$ws.Range("I17").Value = 2999
$ws.Range("K17").Value = 2999
$ws.Range("M17").Value = -2829
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").Value = $null
$ws.Range("H22").Value = 409.8889
$ws.Range("I22").Value = 368.625
$ws.Range("K22").Value = 368.625
$ws.Range("M22").Value = -73.625
$ws.Range("H27").Value = 409.8889
$ws.Range("I27").Value = 368.625
$ws.Range("K27").Value = 368.625
$ws.Range("M27").Value = -261.625
$ws.Range("H61").Value = 5661.75
$ws.Range("I61").Value = 4915.8335
$ws.Range("J61").Value = 7899.5
$ws.Range("K61").Value = 4915.8335
$ws.Range("L61").Value = 7899.5
$ws.Range("M61").Value = -4713.8335
$ws.Range("N61").Value = -8303.5
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null
$ws.Range("H113").Value = 5661.75
$ws.Range("I113").Value = 4915.8335
$ws.Range("J113").Value = 7899.5
$ws.Range("K113").Value = 4915.8335
$ws.Range("L113").Value = 7899.5
$ws.Range("M113").Value = -2745.8335
$ws.Range("N113").Value = -12239.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = $null

